$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.031.21"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.678.03"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.00"
$ws.Range("E5").Value = "  +5.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.05"
$ws.Range("E6").Value = "  +13.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.717"
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.39"
$ws.Range("E10").Value = "  +19.00%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000283"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.34"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.266.41"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.677.65"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.25"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.832.55"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "404.84"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.54"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.28"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  +8.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.04"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.90"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.01"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.35"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.57"
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "46.96"
$ws.Range("E32").Value = "  +10.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.59"
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("E34").Value = "  +6.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "627.88"
$ws.Range("E35").Value = "  +7.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "66.83"
$ws.Range("E36").Value = "  +4.71%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0824"
$ws.Range("E37").Value = "  -6.00%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.408"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.137"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.00"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0439"
$ws.Range("E43").Value = "  +3.05%  "
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.894.01"
$ws.Range("E45").Value = "  +5.40%  "
$ws.Range("E46").Value = "  +5.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.12"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.18"
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.65"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("E51").Value = "  -1.97%  "
